$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("I9").Value = "dhfs"
$ws.Range("I9").Select()
